# feat : 대화 나가기버튼 추가, docs : dialogSheet 수정
# Add two new "leaveTalk" dialog rows (farewell lines) to the DialogData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data rows (13 and 14): leaveTalk / NPC_Roras ---
$ws.Range("A13").Value = "leaveTalk"
$ws.Range("B13").Value = "NPC_Roras"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = "작별"
$ws.Range("E13").Value = "이런 내가 너무 오래 잡아둔건가"

$ws.Range("A14").Value = "leaveTalk"
$ws.Range("B14").Value = "NPC_Roras"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "작별"
$ws.Range("E14").Value = "즐거운 모험 되시게"

# --- Update the list/whole-number data validations so their ranges
#     cover the newly added row 14 as well (and the open-ended trailing
#     ranges now start at row 15 instead of row 14). ---
$ws.Range("A2:A1048576").Validation.Delete()
$ws.Range("A2:A1048576").Validation.Add(3, 1, 1, """dialog,greeting,leaveTalk,quest""")

$ws.Range("C2:C1048576").Validation.Delete()
$ws.Range("C2:C1048576").Validation.Add(1, 1, 1, "-2,147,483,648", "2,147,483,647")

# --- Update the active selection to match the editor's last position ---
$ws.Range("E9").Select()

Write-Output "Added leaveTalk rows 13-14 and refreshed validations/selection"
